# =====================================================================
# "Update gh-pages to output generated at 456a3b4"
#
# A fresh scrape re-lists "广州·第九届娃展沙龙·顽哇展（GZDP9)" as its own
# row (duplicating the row already present for that event) in both the
# "展览" and "全部类型" sheets, pushing every later row down by one. The
# scrape also refreshed several "想去人数" (want-to-go count) numbers
# across all four sheets.
# =====================================================================

$wb = $excel.ActiveWorkbook
$wsExpo = $wb.Worksheets.Item(1)   # 展览
$wsShow = $wb.Worksheets.Item(2)   # 演出
$wsLife = $wb.Worksheets.Item(3)   # 本地生活
$wsAll  = $wb.Worksheets.Item(4)   # 全部类型

function Set-TextValue($range, $text) {
    # Force a literal-text write so digit/date-like strings (e.g. '2024-03-24')
    # are not auto-converted into date serials, matching the inlineStr cells
    # used throughout this workbook.
    $range.Value = "'" + $text
}

function Set-DuplicateRow($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i) {
    Set-TextValue $ws.Range("B$row") $b
    Set-TextValue $ws.Range("C$row") $c
    Set-TextValue $ws.Range("D$row") $d
    Set-TextValue $ws.Range("E$row") $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g
    Set-TextValue $ws.Range("H$row") $h
    Set-TextValue $ws.Range("I$row") $i
}

# --- "展览": insert the new row 14 (duplicates the updated row 13) --------
$wsExpo.Rows.Item(14).Insert()
Set-DuplicateRow $wsExpo 14 "2024-03-24" "广州·第九届娃展沙龙·顽哇展（GZDP9)" "洛浦街厦滘西环路1号 岭南会展中心" "2024.03.24 11:00-03.24 17:00" 171 49.5 "https://show.bilibili.com/platform/detail.html?id=82181" "//i1.hdslb.com/bfs/openplatform/202402/DUqC4OGd1709176548467.jpeg"

# Renumber column A (sequential index, row N -> N-1) through the new last row
for ($r = 14; $r -le 32; $r++) {
    $wsExpo.Cells.Item($r, 1).Value = $r - 1
}

# "想去人数" (F column) refreshes - rows unaffected by the insertion
$wsExpo.Range("F2").Value = 914
$wsExpo.Range("F3").Value = 1480
$wsExpo.Range("F4").Value = 1151
$wsExpo.Range("F5").Value = 538
$wsExpo.Range("F8").Value = 698
$wsExpo.Range("F9").Value = 294
$wsExpo.Range("F13").Value = 171

# "想去人数" (F column) refreshes - rows shifted down by the insertion
$wsExpo.Range("F15").Value = 3864
$wsExpo.Range("F21").Value = 298
$wsExpo.Range("F27").Value = 270
$wsExpo.Range("F30").Value = 1636
$wsExpo.Range("F31").Value = 367

# --- "演出": F column refresh ------------------------------------------------
$wsShow.Range("F3").Value = 35

# --- "本地生活": F column refresh --------------------------------------------
$wsLife.Range("F2").Value = 395

# --- "全部类型": insert the new row 19 (duplicates the updated row 18) -----
$wsAll.Rows.Item(19).Insert()
Set-DuplicateRow $wsAll 19 "2024-03-24" "广州·第九届娃展沙龙·顽哇展（GZDP9)" "洛浦街厦滘西环路1号 岭南会展中心" "2024.03.24 11:00-03.24 17:00" 171 49.5 "https://show.bilibili.com/platform/detail.html?id=82181" "//i1.hdslb.com/bfs/openplatform/202402/DUqC4OGd1709176548467.jpeg"

# Renumber column A (sequential index, row N -> N-1) through the new last row
for ($r = 19; $r -le 45; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

# "想去人数" (F column) refreshes - rows unaffected by the insertion
$wsAll.Range("F2").Value = 395
$wsAll.Range("F3").Value = 914
$wsAll.Range("F4").Value = 1480
$wsAll.Range("F5").Value = 1151
$wsAll.Range("F7").Value = 35
$wsAll.Range("F9").Value = 538
$wsAll.Range("F12").Value = 698
$wsAll.Range("F14").Value = 294
$wsAll.Range("F18").Value = 171

# "想去人数" (F column) refreshes - rows shifted down by the insertion
$wsAll.Range("F20").Value = 3864
$wsAll.Range("F27").Value = 298
$wsAll.Range("F40").Value = 270
$wsAll.Range("F43").Value = 1636
$wsAll.Range("F44").Value = 367

